$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns E and F for rows 1 through 53 (header + 52 data rows)
for ($r = 1; $r -le 53; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2
    $eCell.Value2 = $fVal
    $fCell.Value2 = $eVal
}
